$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Update existing row 8: surname changes from "Demir" to "Canli" (Canlı)
$ws.Range("C8").Value = "Canlı"

# Add new row 9 with a new student record
$ws.Range("A9").Value = 20251006
$ws.Range("B9").Value = "Can"
$ws.Range("C9").Value = "Berk"
$ws.Range("D9").Value = "1B"

# Match formatting of the new row to the rest of the data rows (same as row 8)
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)

$ws.Range("C26").Select()
